$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.278.30'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.73%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.770.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.70%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.17'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.53%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.45'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.05%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.602'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.36%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.80%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.93'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.51%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.66%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.72%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.65%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.46'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.28%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.906.00'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.18%  '

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.75%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.778.78'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.55%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.25'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.15%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.87'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.33%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '361.24'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.19%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.37%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.71%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.533'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -5.76%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.16'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.97%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.65%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.64'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.26%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.11%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0920'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.63%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.44'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.85%  '

$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.39'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +12.45%  '

$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.00'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.28%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '167.67'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.69%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.70%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.03'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.53%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.24'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.83'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.29%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.02'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.18%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '350.46'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.97%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.39'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.14%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.23'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.83%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.23'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.10%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.78'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.17%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.67'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.10%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0597'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.25%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '137.69'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.23%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.38%  '

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.78%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.57%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.43%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.142.55'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.89%  '
